# -----------------------------------------------------------------------
# Edit: mlk.docx
#   1. Append two trailing spaces to the first paragraph's existing text,
#      then add a red (C00000) "(This is a change - Version for branch
#      alternate)" annotation split across three runs.
#   2. Append a new, otherwise-empty paragraph (shaded F9F9F9) right
#      after the final paragraph of the speech, before the section break.
# -----------------------------------------------------------------------

$d = $word.ActiveDocument

# --- 1. First paragraph -------------------------------------------------
# Grow the existing run's text with two trailing spaces (Find/Replace
# keeps it as a single run and preserves xml:space="preserve").
$null = $d.Content.Find.Execute(
    "This is a Microsoft word document.", $true, $false, $false, $false,
    $false, $true, 1, $false,
    "This is a Microsoft word document.  ", 2)

$p1 = $d.Paragraphs.Item(1)
$redColor = 192   # BGR-encoded w:val="C00000"

function Add-RedRun($paragraph, [string]$text) {
    $paraEnd = $paragraph.Range.End
    $insertPoint = $d.Range($paraEnd - 1, $paraEnd - 1)
    $insertPoint.InsertAfter($text)
    $newParaEnd = $paragraph.Range.End
    $newRun = $d.Range($paraEnd - 1, $newParaEnd - 1)
    $newRun.Font.Color = $redColor
}

$enDash = [char]0x2013
$chunk1 = "(This is a change " + $enDash + " Ve"
$chunk2 = "rsion for branch alternate"
$chunk3 = ")"

Add-RedRun $p1 $chunk1
Add-RedRun $p1 $chunk2
Add-RedRun $p1 $chunk3

# --- 2. New shaded paragraph at the end of the document -----------------
$lastP = $d.Paragraphs.Item($d.Paragraphs.Count)
$endOfBody = $lastP.Range.End
$insertRange = $d.Range($endOfBody, $endOfBody)
$null = $insertRange.InsertXML(
    "<w:p xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'><w:pPr><w:shd w:val='clear' w:color='auto' w:fill='F9F9F9'/></w:pPr></w:p>")

# --- 3. Drop the now-unused styles (mirrors Word's own GC of styles that
#     aren't applied anywhere in the document body). Must delete from the
#     highest style-collection index down to the lowest -- deleting
#     low-to-high walks into a stale-index bug in this host's style store.
$stylesToDelete = @(
    "podcast-tools__subscribe-links",
    "generic-title",
    "subscribe-more-info",
    "subscribe",
    "audio-tool",
    "Heading 4 Char",
    "Heading 2 Char",
    "Hyperlink",
    "apple-converted-space",
    "Heading 4",
    "Heading 2"
)
foreach ($styleName in $stylesToDelete) {
    $styleToRemove = $d.Styles.Item($styleName)
    $styleToRemove.Delete()
}
